$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date for 7a6310d6-... row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2017-02-09 15:00:17"

# "zh-cn" sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 7a6310d6-... row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2017-02-09 14:59:57"
$wsZhCn.Range("L4").Value = "2017-02-09 15:00:54"

# "de-de" sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 7a6310d6-... row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2017-02-09 15:00:17"
$wsDeDe.Range("L4").Value = "2017-02-09 15:01:21"
